# Added team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should use the same style as the rest of row 1 (bold,
# centered, bordered). Copy that formatting from the last existing header
# cell (AB1) onto the three new header cells before setting their text.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every player row (2-37) gets the same team record: 49 wins, 64 losses,
# 0 ties.
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 29).Value = 49
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 0
}
